# MassWateR_Results_Template.xlsx update
#
# The "Instructions" worksheet documents each column of the "Results" sheet.
# Two of the instruction descriptions were expanded with additional guidance:
#
#   1. "Characteristic Name" row (B13) - clarified that either the WQX or
#      Simple parameter name can be used.
#   2. "Result Attribute" row (B19) - clarified that the field is optional
#      "(for analysis)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

$ws.Range("B13").Value = "Name of the measured parameter.`r`n- Note that this can be either the WQX or Simple parameter name.  However, if a parameter is distinguished by Sample Fraction only (i.e. TDP, PON, TDN), then the Simple parameter name must be used here and in all other files (WQX, DQO, etc.)."

$ws.Range("B19").Value = "Optional (for analysis) field that can be used to group results using whatever coding system the users wants, such as wet weather or dry weather.`r`nThis field can also be used to record sample IDs for QC samples, such as lab blanks and lab dulicates."
